$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds two API-spec blocks:
#   rows 1-9   : getCuTrusByName(hoTen)
#   rows 10-18 : getCuTruById(id)
# We append a third block (rows 19-27) for the new getExpiredCuTrus(hetHan)
# API. It has exactly the same shape/styling as the getCuTruById block
# (rows 10-18), so we copy it cell-by-cell (this preserves the existing
# cell styles exactly instead of generating new ones) and then only change
# the two cells whose text actually differs (B19 and B20).

for ($i = 0; $i -lt 9; $i++) {
    $srcRow = 10 + $i
    $dstRow = 19 + $i

    $srcA = $ws.Cells.Item($srcRow, 1)
    $dstA = $ws.Cells.Item($dstRow, 1)
    $srcA.Copy($dstA)

    $srcB = $ws.Cells.Item($srcRow, 2)
    $dstB = $ws.Cells.Item($dstRow, 2)
    $srcB.Copy($dstB)
}

# Row 19 / col B: API name
$ws.Cells.Item(19, 2).Value = "getExpiredCuTrus(hetHan)"

# Row 20 / col B: description with the italicized param name "hetHan"
$descCell = $ws.Cells.Item(20, 2)
$descCell.Value = "Lấy ra object CuTru dựa trên param hetHan"
$paramChars = $descCell.Characters(36, 6)
$paramChars.Font.Italic = $true
$paramChars.Font.Name = "Times New Roman"
$paramChars.Font.Size = 14

# Re-create the merged header cell for the new block (A20:A25), mirroring
# the existing A11:A16 merge used by the block we copied from.
$ws.Range("A20:A25").Merge()

# Update the window view to match the authored state: scrolled so row 11
# is at the top, with B22 selected.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("B22").Select() | Out-Null
